# SwaadSutra_Consolidated_2026-01-20.xlsx update
# A new order (#25, Indrani Karmakar) came in on 2026-01-20 at 15:05.
# It is inserted as the new top data row (row 2) of "All Orders",
# pushing every existing order down by one row, and the "Daily Summary"
# totals for 2026-01-20 are bumped to reflect the new order.

$wb = $excel.ActiveWorkbook

# ---- All Orders sheet -----------------------------------------------
$ws = $wb.Worksheets.Item("All Orders")

# Push all existing order rows down by one to make room for the newest order.
$ws.Rows(2).Insert()

# Helper: write a value as genuine text (even when it looks numeric),
# matching this workbook's convention of storing Order ID lookalikes
# (phone numbers, ISO dates, etc.) as text, then drop back to the
# workbook's default "Normal" style so no new per-cell number format
# is left behind.
function Set-TextValue($cell, $text) {
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

$row = $ws.Cells.Item(2, 1)
$row.Value = 25

Set-TextValue $ws.Cells.Item(2, 2) "2026-01-20 15:05"
Set-TextValue $ws.Cells.Item(2, 3) "Indrani Karmakar"
Set-TextValue $ws.Cells.Item(2, 4) "A-1903"
Set-TextValue $ws.Cells.Item(2, 5) "7030961520"
Set-TextValue $ws.Cells.Item(2, 6) "Appe Chutney x1"

$ws.Cells.Item(2, 7).Value = 60

Set-TextValue $ws.Cells.Item(2, 8) "NEW"
Set-TextValue $ws.Cells.Item(2, 9) "PENDING"
Set-TextValue $ws.Cells.Item(2, 10) "2026-01-21"
Set-TextValue $ws.Cells.Item(2, 11) "09:00"
Set-TextValue $ws.Cells.Item(2, 12) "Less spicy. Flavourful"
Set-TextValue $ws.Cells.Item(2, 13) ""
Set-TextValue $ws.Cells.Item(2, 14) ""

# ---- Daily Summary sheet ---------------------------------------------
# Reflect the new order in the 2026-01-20 roll-up: one more total order,
# and its 60 (still PENDING) add to both Revenue and Pending.
$ws2 = $wb.Worksheets.Item("Daily Summary")

$ws2.Cells.Item(2, 2).Value = 6    # Total Orders: 5 -> 6
$ws2.Cells.Item(2, 5).Value = 380  # Revenue: 320 -> 380
$ws2.Cells.Item(2, 7).Value = 330  # Pending: 270 -> 330
